$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.525.25'
$ws.Range("E2").Value = '  +0.54%  '

$ws.Range("D3").Value = '3.146.46'
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'" + '610.27'
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("D6").Value = "'" + '144.28'
$ws.Range("E6").Value = '  -2.21%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.142.75'
$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").Value = "'" + '0.526'
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("D11").Value = "'" + '5.41'
$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").Value = "'" + '0.468'
$ws.Range("E12").Value = '  -1.23%  '

$ws.Range("D13").Value = "'" + '0.0000259'
$ws.Range("E13").Value = '  +2.99%  '

$ws.Range("D14").Value = "'" + '35.50'
$ws.Range("E14").Value = '  -0.14%  '

$ws.Range("D15").Value = '3.657.72'
$ws.Range("E15").Value = '  +0.13%  '

$ws.Range("E16").Value = '  +2.45%  '

$ws.Range("D17").Value = '64.362.65'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '3.143.59'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").Value = "'" + '6.87'
$ws.Range("E19").Value = '  -0.42%  '

$ws.Range("D20").Value = "'" + '481.19'
$ws.Range("E20").Value = '  +0.51%  '

$ws.Range("D21").Value = "'" + '14.74'
$ws.Range("E21").Value = '  -0.56%  '

$ws.Range("D22").Value = "'" + '0.718'

$ws.Range("D23").Value = "'" + '7.73'
$ws.Range("E23").Value = '  -1.19%  '

$ws.Range("D24").Value = "'" + '85.54'
$ws.Range("E24").Value = '  +2.96%  '

$ws.Range("D25").Value = "'" + '13.42'
$ws.Range("E25").Value = '  -1.91%  '

$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  -3.97%  '

$ws.Range("D28").Value = "'" + '8.48'
$ws.Range("E28").Value = '  -0.14%  '

$ws.Range("D29").Value = "'" + '7.21'
$ws.Range("E29").Value = '  +5.93%  '

$ws.Range("D30").Value = "'" + '0.115'
$ws.Range("E30").Value = '  +1.31%  '

$ws.Range("E31").Value = '  -6.33%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'" + '27.00'
$ws.Range("E32").Value = '  +2.87%  '

$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = "'" + '1.00'
$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("D34").Value = "'" + '2.66'
$ws.Range("E34").Value = '  -3.26%  '

$ws.Range("E35").Value = '  -0.33%  '

$ws.Range("D36").Value = "'" + '5.99'
$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("D37").Value = '0.0₃0770'
$ws.Range("E37").Value = '  +4.24%  '

$ws.Range("D38").Value = "'" + '52.53'

$ws.Range("D39").Value = "'" + '3.05'
$ws.Range("E39").Value = '  +3.42%  '

$ws.Range("D40").Value = "'" + '447.19'
$ws.Range("E40").Value = '  -2.70%  '

$ws.Range("D41").Value = "'" + '0.0393'
$ws.Range("E41").Value = '  -0.60%  '

$ws.Range("D42").Value = "'" + '0.120'
$ws.Range("E42").Value = '  +1.16%  '

$ws.Range("D43").Value = "'" + '8.26'
$ws.Range("E43").Value = '  -2.08%  '

$ws.Range("D44").Value = '2.873.36'
$ws.Range("E44").Value = '  +0.85%  '

$ws.Range("D45").Value = "'" + '0.262'
$ws.Range("E45").Value = '  -1.43%  '

$ws.Range("E46").Value = '  -1.20%  '

$ws.Range("D47").Value = "'" + '2.42'
$ws.Range("E47").Value = '  +2.62%  '

$ws.Range("E48").Value = '  +0.06%  '

$ws.Range("D49").Value = "'" + '26.24'
$ws.Range("E49").Value = '  -1.12%  '

$ws.Range("E50").Value = '  -0.81%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = "'" + '119.81'
$ws.Range("E51").Value = '  +0.24%  '
